# Indonesia Liga 1 - base update (08-05-2024 20:15)
# The source feed re-ordered two pairs of fixtures that share the same
# match date/time; the row index (column A) and everything else about the
# sheet stays put, only the match data in columns B:AB is swapped between
# the two rows in each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AB$rowA")
    $rangeB = $ws.Range("B$rowB`:AB$rowB")

    $valA = $rangeA.Value2
    $valB = $rangeB.Value2

    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# Rows 300/301 (ids 298/299) swap their fixture data.
Swap-RowData 300 301

# Rows 304/305 (ids 302/303) swap their fixture data.
Swap-RowData 304 305
